$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: TDRIssuanceMaturityLCY_IBG
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(1,6).Value = "DRAWDOWN.ACCOUNT"
$ws1.Cells.Item(1,7).Value = "PRIN.LIQ.ACCT"
$ws1.Cells.Item(1,8).Value = "INT.LIQ.ACCT"

$ws1.Range("F2:H2").NumberFormat = "@"
$ws1.Cells.Item(2,6).Value = "5000000522"
$ws1.Cells.Item(2,7).Value = "5000000522"
$ws1.Cells.Item(2,8).Value = "5000000522"

$ws1.Columns.Item(6).ColumnWidth = 21.43
$ws1.Columns.Item(7).ColumnWidth = 13.1
$ws1.Columns.Item(8).ColumnWidth = 11.6

$ws1.Range("F1:H2").Select()

# ---------------------------------------------------------------------------
# Sheet 2: TDRIssuanceMaturityFCY_IBG
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2,1).Value = 11871146

$ws2.Cells.Item(1,6).Value = "DRAWDOWN.ACCOUNT"
$ws2.Cells.Item(1,7).Value = "PRIN.LIQ.ACCT"
$ws2.Cells.Item(1,8).Value = "INT.LIQ.ACCT"

$ws2.Range("F2:H2").NumberFormat = "@"
$ws2.Cells.Item(2,6).Value = "5000000529"
$ws2.Cells.Item(2,7).Value = "5000000529"
$ws2.Cells.Item(2,8).Value = "5000000529"

$ws2.Columns.Item(6).ColumnWidth = 21.43
$ws2.Columns.Item(7).ColumnWidth = 13.1
$ws2.Columns.Item(8).ColumnWidth = 11.6

$ws2.Range("F2").Select()

# ---------------------------------------------------------------------------
# Sheet 3: TDRBackDatedMaturityLCY_IBG
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2,1).Value = 14337500

$ws3.Cells.Item(1,7).Value = "DRAWDOWN.ACCOUNT"
$ws3.Cells.Item(1,8).Value = "PRIN.LIQ.ACCT"
$ws3.Cells.Item(1,9).Value = "INT.LIQ.ACCT"

$ws3.Range("G2:I2").NumberFormat = "@"
$ws3.Cells.Item(2,7).Value = "5000000531"
$ws3.Cells.Item(2,8).Value = "5000000531"
$ws3.Cells.Item(2,9).Value = "5000000531"

$ws3.Columns.Item(7).ColumnWidth = 21.43
$ws3.Columns.Item(8).ColumnWidth = 13.1
$ws3.Columns.Item(9).ColumnWidth = 11.6

$ws3.Range("G1:I2").Select()
$ws3.Cells.Item(2,7).Select()

# ---------------------------------------------------------------------------
# Sheet 4: TDRBackDatedMaturityFCY_IBG
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2,1).Value = 11871196

$ws4.Cells.Item(1,7).Value = "DRAWDOWN.ACCOUNT"
$ws4.Cells.Item(1,8).Value = "PRIN.LIQ.ACCT"
$ws4.Cells.Item(1,9).Value = "INT.LIQ.ACCT"

$ws4.Range("G2:I2").NumberFormat = "@"
$ws4.Cells.Item(2,7).Value = "5000000534"
$ws4.Cells.Item(2,8).Value = "5000000534"
$ws4.Cells.Item(2,9).Value = "5000000534"

$ws4.Columns.Item(7).ColumnWidth = 21.43
$ws4.Columns.Item(8).ColumnWidth = 13.1
$ws4.Columns.Item(9).ColumnWidth = 11.6

$ws4.Activate()
$ws4.Range("G1:I2").Select()
$ws4.Cells.Item(1,7).Select()
